# Auto-generated edit script applying the Atomos_Profits.xlsx diff
# Updates numeric value cells (H..N) across multiple leve-profit worksheets.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H15").Value = 1999
$ws_ALC.Range("I15").Value = 1999
$ws_ALC.Range("K15").Value = 5997
$ws_ALC.Range("M15").Value = -5828
$ws_ALC.Range("H127").Value = 1699756.6
$ws_ALC.Range("I127").Value = 2835.375
$ws_ALC.Range("J127").Value = 2933881.2
$ws_ALC.Range("K127").Value = 8506.125
$ws_ALC.Range("L127").Value = 8801643.600000001
$ws_ALC.Range("M127").Value = -3546.125
$ws_ALC.Range("N127").Value = -8811563.600000001
$ws_ALC.Range("H129").Value = 992.6923
$ws_ALC.Range("J129").Value = 889.5278
$ws_ALC.Range("L129").Value = 2668.5834
$ws_ALC.Range("N129").Value = -12668.5834
$ws_ALC.Range("H135").Value = 955.8929000000001
$ws_ALC.Range("I135").Value = 508.6842
$ws_ALC.Range("J135").Value = 1900
$ws_ALC.Range("K135").Value = 4578.1578
$ws_ALC.Range("L135").Value = 17100
$ws_ALC.Range("M135").Value = -2043.1578
$ws_ALC.Range("N135").Value = -22170
$ws_ALC.Range("H138").Value = 5968.3945
$ws_ALC.Range("I138").Value = 2741.9644
$ws_ALC.Range("J138").Value = 7850.479
$ws_ALC.Range("K138").Value = 8225.893199999999
$ws_ALC.Range("L138").Value = 23551.437
$ws_ALC.Range("M138").Value = -3085.893199999999
$ws_ALC.Range("N138").Value = -33831.43700000001

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H23").Value = 0
$ws_ARM.Range("J23").Value = 0
$ws_ARM.Range("L23").Value = 0
$ws_ARM.Range("N23").ClearContents()
$ws_ARM.Range("H32").Value = 19260.75
$ws_ARM.Range("I32").Value = 15553.718
$ws_ARM.Range("J32").Value = 30381.846
$ws_ARM.Range("K32").Value = 15553.718
$ws_ARM.Range("L32").Value = 30381.846
$ws_ARM.Range("M32").Value = -15266.718
$ws_ARM.Range("N32").Value = -30955.846
$ws_ARM.Range("H37").Value = 18038
$ws_ARM.Range("J37").Value = 18038
$ws_ARM.Range("L37").Value = 18038
$ws_ARM.Range("N37").Value = -18584
$ws_ARM.Range("H61").Value = 2771.5532
$ws_ARM.Range("I61").Value = 2099.3333
$ws_ARM.Range("J61").Value = 3679.05
$ws_ARM.Range("K61").Value = 2099.3333
$ws_ARM.Range("L61").Value = 3679.05
$ws_ARM.Range("M61").Value = -1887.3333
$ws_ARM.Range("N61").Value = -4103.05
$ws_ARM.Range("H74").Value = 1752.2307
$ws_ARM.Range("I74").Value = 1257.2273
$ws_ARM.Range("J74").Value = 4474.75
$ws_ARM.Range("K74").Value = 1257.2273
$ws_ARM.Range("L74").Value = 4474.75
$ws_ARM.Range("M74").Value = -383.2273
$ws_ARM.Range("N74").Value = -6222.75
$ws_ARM.Range("H77").Value = 1752.2307
$ws_ARM.Range("I77").Value = 1257.2273
$ws_ARM.Range("J77").Value = 4474.75
$ws_ARM.Range("K77").Value = 6286.136500000001
$ws_ARM.Range("L77").Value = 22373.75
$ws_ARM.Range("M77").Value = -1918.136500000001
$ws_ARM.Range("N77").Value = -31109.75
$ws_ARM.Range("H132").Value = 3121.7585
$ws_ARM.Range("I132").Value = 2746.7693
$ws_ARM.Range("J132").Value = 3426.4375
$ws_ARM.Range("K132").Value = 8240.3079
$ws_ARM.Range("L132").Value = 10279.3125
$ws_ARM.Range("M132").Value = -5710.3079
$ws_ARM.Range("N132").Value = -15339.3125
$ws_ARM.Range("H136").Value = 2771.5532
$ws_ARM.Range("I136").Value = 2099.3333
$ws_ARM.Range("J136").Value = 3679.05
$ws_ARM.Range("K136").Value = 6297.999899999999
$ws_ARM.Range("L136").Value = 11037.15
$ws_ARM.Range("M136").Value = -3747.999899999999
$ws_ARM.Range("N136").Value = -16137.15

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H82").Value = 23468.309
$ws_BSM.Range("J82").Value = 29028.9
$ws_BSM.Range("L82").Value = 29028.9
$ws_BSM.Range("N82").Value = -29794.9
$ws_BSM.Range("H85").Value = 23468.309
$ws_BSM.Range("J85").Value = 29028.9
$ws_BSM.Range("L85").Value = 29028.9
$ws_BSM.Range("N85").Value = -31680.9
$ws_BSM.Range("H134").Value = 3085.8958
$ws_BSM.Range("I134").Value = 2779.3489
$ws_BSM.Range("J134").Value = 5722.2
$ws_BSM.Range("K134").Value = 8338.046699999999
$ws_BSM.Range("L134").Value = 17166.6
$ws_BSM.Range("M134").Value = -5803.046699999999
$ws_BSM.Range("N134").Value = -22236.6

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H22").Value = 1118
$ws_CRP.Range("I22").Value = 482.63635
$ws_CRP.Range("J22").Value = 2515.8
$ws_CRP.Range("K22").Value = 482.63635
$ws_CRP.Range("L22").Value = 2515.8
$ws_CRP.Range("M22").Value = -132.63635
$ws_CRP.Range("N22").Value = -3215.8
$ws_CRP.Range("H31").Value = 1616009.2
$ws_CRP.Range("I31").Value = 2382733.5
$ws_CRP.Range("J31").Value = 5888.4
$ws_CRP.Range("K31").Value = 2382733.5
$ws_CRP.Range("L31").Value = 5888.4
$ws_CRP.Range("M31").Value = -2382438.5
$ws_CRP.Range("N31").Value = -6478.4
$ws_CRP.Range("H34").Value = 1616009.2
$ws_CRP.Range("I34").Value = 2382733.5
$ws_CRP.Range("J34").Value = 5888.4
$ws_CRP.Range("K34").Value = 2382733.5
$ws_CRP.Range("L34").Value = 5888.4
$ws_CRP.Range("M34").Value = -2382531.5
$ws_CRP.Range("N34").Value = -6292.4
$ws_CRP.Range("H141").Value = 29792.105
$ws_CRP.Range("J141").Value = 32011.428
$ws_CRP.Range("L141").Value = 32011.428
$ws_CRP.Range("N141").Value = -42371.428

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 795.5961
$ws_CUL.Range("I5").Value = 503.5625
$ws_CUL.Range("J5").Value = 4300
$ws_CUL.Range("K5").Value = 1510.6875
$ws_CUL.Range("L5").Value = 12900
$ws_CUL.Range("M5").Value = -1398.6875
$ws_CUL.Range("N5").Value = -13124
$ws_CUL.Range("H114").Value = 934.875
$ws_CUL.Range("I114").Value = 353.5
$ws_CUL.Range("J114").Value = 1516.25
$ws_CUL.Range("K114").Value = 1060.5
$ws_CUL.Range("L114").Value = 4548.75
$ws_CUL.Range("M114").Value = 2193.5
$ws_CUL.Range("N114").Value = -11056.75
$ws_CUL.Range("H131").Value = 1327.8857
$ws_CUL.Range("J131").Value = 1105.8524
$ws_CUL.Range("L131").Value = 3317.5572
$ws_CUL.Range("N131").Value = -13397.5572
$ws_CUL.Range("H135").Value = 795.5961
$ws_CUL.Range("I135").Value = 503.5625
$ws_CUL.Range("J135").Value = 4300
$ws_CUL.Range("K135").Value = 4532.0625
$ws_CUL.Range("L135").Value = 38700
$ws_CUL.Range("M135").Value = -1997.0625
$ws_CUL.Range("N135").Value = -43770

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H92").Value = 9579.799999999999
$ws_GSM.Range("J92").Value = 9579.799999999999
$ws_GSM.Range("L92").Value = 9579.799999999999
$ws_GSM.Range("N92").Value = -13323.8
$ws_GSM.Range("H132").Value = 4318.6055
$ws_GSM.Range("I132").Value = 5191.75
$ws_GSM.Range("J132").Value = 3348.4443
$ws_GSM.Range("K132").Value = 15575.25
$ws_GSM.Range("L132").Value = 10045.3329
$ws_GSM.Range("M132").Value = -13045.25
$ws_GSM.Range("N132").Value = -15105.3329
$ws_GSM.Range("H136").Value = 28157
$ws_GSM.Range("J136").Value = 28157
$ws_GSM.Range("L136").Value = 84471
$ws_GSM.Range("N136").Value = -89571

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H55").Value = 1138.6428
$ws_LTW.Range("I55").Value = 254.1
$ws_LTW.Range("J55").Value = 3350
$ws_LTW.Range("K55").Value = 254.1
$ws_LTW.Range("L55").Value = 3350
$ws_LTW.Range("M55").Value = -81.09999999999999
$ws_LTW.Range("N55").Value = -3696

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H132").Value = 2338994.2
$ws_WVR.Range("I132").Value = 2781737.2
$ws_WVR.Range("J132").Value = 62029.855
$ws_WVR.Range("K132").Value = 8345211.600000001
$ws_WVR.Range("L132").Value = 186089.565
$ws_WVR.Range("M132").Value = -8342681.600000001
$ws_WVR.Range("N132").Value = -191149.565
$ws_WVR.Range("H138").Value = 29571.75
$ws_WVR.Range("J138").Value = 29571.75
$ws_WVR.Range("L138").Value = 29571.75
$ws_WVR.Range("N138").Value = -39851.75
